# Add a new "ExcelId" identifier column at the front of the table (shifting
# the existing ExcelProductName/ExcelPrice/ExcelDeliveryDays/ExcelDescription/
# ExcelDiscount columns from A:E to B:F), populate the new column's values,
# and update the first product's name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing B..F data/styles shift right
# automatically and keep their original formatting.
$ws.Range("A1").EntireColumn.Insert()

# Match the new column's look (alignment / width) to the rest of the table.
$ws.Range("A1:A4").HorizontalAlignment = 1
$ws.Columns.Item(1).ColumnWidth = 12.74

# Header
$ws.Range("A1").Value = "ExcelId"

# Row 2 - "Widget Excel A" -> gets a real id and a renamed product name
$ws.Range("A2").Value = "d0a8bec1-2690-437c-b857-3a399df25b83"
$ws.Range("B2").Value = "Widget Excel A update"

# Row 3 - "Widget Excel B" -> id is a single space placeholder
$ws.Range("A3").Value = " "

# Row 4 - "Widget Excel C" -> id left empty
$ws.Range("A4").Value = ""
